$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2880.4167
$ws.Range("J111").Value = 3885.6667
$ws.Range("L111").Value = 11657.0001
$ws.Range("N111").Value = -17791.0001

$ws.Range("H125").Value = 1558.4
$ws.Range("I125").Value = 1271.909
$ws.Range("K125").Value = 11447.181
$ws.Range("M125").Value = -8987.181

$ws.Range("H132").Value = 1517.5581
$ws.Range("I132").Value = 1217.3948
$ws.Range("K132").Value = 3652.1844
$ws.Range("M132").Value = -1122.1844

$ws.Range("H138").Value = 1084791
$ws.Range("I138").Value = 555.087
$ws.Range("J138").Value = 1493601.2
$ws.Range("K138").Value = 1665.261
$ws.Range("L138").Value = 4480803.6
$ws.Range("M138").Value = 3474.739
$ws.Range("N138").Value = -4491083.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 915.8
$ws.Range("I97").Value = 928.1667
$ws.Range("K97").Value = 928.1667
$ws.Range("M97").Value = -432.1667

$ws.Range("H122").Value = 1767.9546
$ws.Range("I122").Value = 1468.2106
$ws.Range("K122").Value = 4404.6318
$ws.Range("M122").Value = -1954.6318

$ws.Range("H132").Value = 6548.88
$ws.Range("I132").Value = 3525.611
$ws.Range("K132").Value = 10576.833
$ws.Range("M132").Value = -8046.832999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3422.389
$ws.Range("I86").Value = 3125.3635
$ws.Range("J86").Value = 3889.1428
$ws.Range("K86").Value = 3125.3635
$ws.Range("L86").Value = 3889.1428
$ws.Range("M86").Value = -2002.3635
$ws.Range("N86").Value = -6135.1428

$ws.Range("H89").Value = 3422.389
$ws.Range("I89").Value = 3125.3635
$ws.Range("J89").Value = 3889.1428
$ws.Range("K89").Value = 15626.8175
$ws.Range("L89").Value = 19445.714
$ws.Range("M89").Value = -10010.8175
$ws.Range("N89").Value = -30677.714

$ws.Range("H99").Value = 2450.1428
$ws.Range("I99").Value = 1860.2
$ws.Range("J99").Value = 3925
$ws.Range("K99").Value = 1860.2
$ws.Range("L99").Value = 3925
$ws.Range("M99").Value = -362.2
$ws.Range("N99").Value = -6921

$ws.Range("H134").Value = 33786.13
$ws.Range("I134").Value = 1444.4814
$ws.Range("K134").Value = 4333.4442
$ws.Range("M134").Value = -1798.4442

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 299.5
$ws.Range("I22").Value = 299.5
$ws.Range("K22").Value = 299.5
$ws.Range("M22").Value = 50.5

$ws.Range("H31").Value = 835062.9399999999
$ws.Range("I31").Value = 1408.3334
$ws.Range("K31").Value = 1408.3334
$ws.Range("M31").Value = -1113.3334

$ws.Range("H34").Value = 835062.9399999999
$ws.Range("I34").Value = 1408.3334
$ws.Range("K34").Value = 1408.3334
$ws.Range("M34").Value = -1206.3334

$ws.Range("H132").Value = 4955.25
$ws.Range("I132").Value = 4955.25
$ws.Range("K132").Value = 14865.75
$ws.Range("M132").Value = -12335.75

$ws.Range("H134").Value = 671131.8
$ws.Range("I134").Value = 910815.4399999999
$ws.Range("J134").Value = 12001.75
$ws.Range("K134").Value = 2732446.32
$ws.Range("L134").Value = 36005.25
$ws.Range("M134").Value = -2729911.32
$ws.Range("N134").Value = -41075.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3307.3635
$ws.Range("I68").Value = 3200
$ws.Range("J68").Value = 3368.7144
$ws.Range("K68").Value = 9600
$ws.Range("L68").Value = 10106.1432
$ws.Range("M68").Value = -8789
$ws.Range("N68").Value = -11728.1432

$ws.Range("H71").Value = 3307.3635
$ws.Range("I71").Value = 3200
$ws.Range("J71").Value = 3368.7144
$ws.Range("K71").Value = 28800
$ws.Range("L71").Value = 30318.4296
$ws.Range("M71").Value = -24744
$ws.Range("N71").Value = -38430.4296

$ws.Range("H98").Value = 785.8
$ws.Range("J98").Value = 785.8
$ws.Range("L98").Value = 2357.4
$ws.Range("N98").Value = -5353.4

$ws.Range("H114").Value = 2730.5
$ws.Range("J114").Value = 2749
$ws.Range("L114").Value = 8247
$ws.Range("N114").Value = -14755

$ws.Range("H122").Value = 1066.45
$ws.Range("I122").Value = 706.8333
$ws.Range("J122").Value = 1220.5714
$ws.Range("K122").Value = 6361.4997
$ws.Range("L122").Value = 10985.1426
$ws.Range("M122").Value = -3911.4997
$ws.Range("N122").Value = -15885.1426

$ws.Range("H134").Value = 3852.5
$ws.Range("I134").Value = 2516.25
$ws.Range("K134").Value = 7548.75
$ws.Range("M134").Value = -2478.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6669
$ws.Range("I70").Value = 5004
$ws.Range("J70").Value = 9999
$ws.Range("K70").Value = 5004
$ws.Range("L70").Value = 9999
$ws.Range("M70").Value = -4734
$ws.Range("N70").Value = -10539

$ws.Range("H73").Value = 6669
$ws.Range("I73").Value = 5004
$ws.Range("J73").Value = 9999
$ws.Range("K73").Value = 5004
$ws.Range("L73").Value = 9999
$ws.Range("M73").Value = -4068
$ws.Range("N73").Value = -11871

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2192.3572
$ws.Range("I16").Value = 1754.4546
$ws.Range("K16").Value = 1754.4546
$ws.Range("M16").Value = -1584.4546

$ws.Range("H46").Value = 3125.742
$ws.Range("I46").Value = 2456
$ws.Range("K46").Value = 2456
$ws.Range("M46").Value = -2268

$ws.Range("H55").Value = 71428920
$ws.Range("I55").Value = 83333690
$ws.Range("J55").Value = 316.5
$ws.Range("K55").Value = 83333690
$ws.Range("L55").Value = 316.5
$ws.Range("M55").Value = -83333517
$ws.Range("N55").Value = -662.5

$ws.Range("H82").Value = 699.5454999999999
$ws.Range("I82").Value = 599.75
$ws.Range("K82").Value = 599.75
$ws.Range("M82").Value = -238.75

$ws.Range("H85").Value = 699.5454999999999
$ws.Range("I85").Value = 599.75
$ws.Range("K85").Value = 599.75
$ws.Range("M85").Value = 648.25

$ws.Range("H115").Value = 85629.664
$ws.Range("J115").Value = 85629.664
$ws.Range("L115").Value = 85629.664
$ws.Range("N115").Value = -87979.664

$ws.Range("H122").Value = 5929.2085
$ws.Range("I122").Value = 5183.3887
$ws.Range("J122").Value = 8166.6665
$ws.Range("K122").Value = 15550.1661
$ws.Range("L122").Value = 24499.9995
$ws.Range("M122").Value = -13100.1661
$ws.Range("N122").Value = -29399.9995

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 119000
$ws.Range("J110").Value = 119000
$ws.Range("L110").Value = 119000
$ws.Range("N110").Value = -127180
